$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename cell-type labels in column A (Level1Classification) to their
# updated wording.
$ws.Range("A17").Value = "Myeloid cells (other)"
$ws.Range("A19").Value = "Leukocytes (other)"
$ws.Range("A29").Value = "Progenitor immune cells"
$ws.Range("A32").Value = "Unclassified cells"
$ws.Range("A36").Value = "Oligodendrocytes"

# Move the active selection to A34, matching the saved view state.
$ws.Range("A34").Select()
